# "Logged Week 16 and performed season sim from Week 17"
# Updates Rushing + Receiving tables with new weekly totals and re-sorts
# each table into its new row order (players matched by name, not by
# previous row position).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Rushing sheet
# ---------------------------------------------------------------------
$rushing = $wb.Worksheets.Item("Rushing")

# name, 1DATT, 2DATT, 3DATT, RZATT
$rushingRows = @(
    @("R.Wilson",   10, 9,  13, 2),
    @("R.Penny",    41, 19, 5,  8),
    @("A.Collins",  55, 36, 9,  12),
    @("D.Dallas",   3,  3,  4,  2),
    @("A.Peterson", 7,  2,  1,  5),
    @("T.Lockett",  1,  0,  0,  0),
    @("F.Swain",    3,  2,  0,  0),
    @("D.Eskridge", 4,  0,  0,  0),
    @("G.Everett",  2,  1,  0,  0),
    @("W.Dissly",   0,  1,  0,  0)
)

$r = 2
foreach ($row in $rushingRows) {
    $rushing.Cells.Item($r, 1).Value = ($r - 2)
    $rushing.Cells.Item($r, 2).Value = $row[0]
    $rushing.Cells.Item($r, 3).Value = $row[1]
    $rushing.Cells.Item($r, 4).Value = $row[2]
    $rushing.Cells.Item($r, 5).Value = $row[3]
    $rushing.Cells.Item($r, 6).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# Receiving sheet
# ---------------------------------------------------------------------
$receiving = $wb.Worksheets.Item("Receiving")

# name, ShortTarget, ShortComp, DeepTarget, DeepComp, RZTarget, RZComp
$receivingRows = @(
    @("R.Penny",     11, 9,  1,  0,  1,  1),
    @("A.Collins",   16, 14, 1,  1,  0,  0),
    @("D.Dallas",    13, 11, 0,  0,  4,  3),
    @("N.Bellore",   1,  1,  0,  0,  0,  0),
    @("D.Metcalf",   76, 54, 33, 10, 15, 10),
    @("T.Lockett",   65, 49, 43, 24, 7,  3),
    @("F.Swain",     30, 18, 7,  4,  4,  1),
    @("D.Eskridge",  11, 9,  6,  0,  4,  2),
    @("P.Hart",      9,  6,  2,  0,  1,  0),
    @("G.Everett",   51, 40, 4,  4,  7,  4),
    @("W.Dissly",    20, 15, 5,  5,  3,  2),
    @("C.Parkinson", 3,  2,  1,  0,  1,  0)
)

$r = 2
foreach ($row in $receivingRows) {
    $receiving.Cells.Item($r, 1).Value = ($r - 2)
    $receiving.Cells.Item($r, 2).Value = $row[0]
    $receiving.Cells.Item($r, 3).Value = $row[1]
    $receiving.Cells.Item($r, 4).Value = $row[2]
    $receiving.Cells.Item($r, 5).Value = $row[3]
    $receiving.Cells.Item($r, 6).Value = $row[4]
    $receiving.Cells.Item($r, 7).Value = $row[5]
    $receiving.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# Last row (C.Parkinson, row 13) picks up the same border formatting as
# the rest of column A instead of its old "no top/bottom border" style.
$receiving.Range("A12").Copy()
[void]$receiving.Range("A13").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Re-activate the Rushing tab (was Receiving) and reset selections to A1
# ---------------------------------------------------------------------
[void]$receiving.Range("A1").Select()
$rushing.Activate()
[void]$rushing.Range("A1").Select()
